$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Details")

# Fill in row 4 with the new test user data (InvalidEmailApp test case)
$ws.Range("A4").Value = "Test User"
$ws.Range("B4").Value = "Mctestuser92.gmail.com"
$ws.Range("C4").Value = "demo123user"
$ws.Range("D4").Value = "demo123user"

# Update the selection to B4 (matches the post-edit cursor position)
$ws.Range("B4").Select()
